$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.510.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.196.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.57%  '
$ws.Range("E6").Value = '  +1.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.83'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.02%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.587'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0912'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.82'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.23%  '
$ws.Range("E13").Value = '  +1.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.529.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.199.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.776'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.451.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000103'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("E21").Value = '  +2.76%  '
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '228.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E27").Value = '  +1.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0803'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.43%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("E36").Value = '  +2.39%  '
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0333'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.08'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("E41").Value = '  +4.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.16%  '
$ws.Range("B45").Value = 'WOONetwork'
$ws.Range("C45").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.474'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +21.34%  '
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0975'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.99%  '
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("E51").Value = '  +1.27%  '
